$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# "Good Morning" -> "GIT UPDATE" for the R10 row's Greeting cell.
$ws.Range("E8").Value = "GIT UPDATE"

# Leave the cursor on the edited cell, matching the saved selection state.
$ws.Range("E8").Select()
